$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Update column C ("Förändrad") from row 2 to row 456 from 45192 to 45202
for ($r = 2; $r -le 456; $r++) {
    $cell = $ws.Cells.Item($r, 3)
    if ($cell.Value2 -eq 45192) {
        $cell.Value = 45202
    }
}

# 2. Row 456 gains an explicit row height (ht="15" customHeight="1")
$ws.Rows.Item(456).RowHeight = 15

# 3. Add new row 457
$ws.Cells.Item(457, 1).Value = "A 45987-2023"
$ws.Cells.Item(457, 2).Value = 45196
$ws.Cells.Item(457, 2).NumberFormat = "YYYY-MM-DD"
$ws.Cells.Item(457, 3).Value = 45202
$ws.Cells.Item(457, 3).NumberFormat = "YYYY-MM-DD"
$ws.Cells.Item(457, 4).Value = "VÄRMLANDS LÄN"
$ws.Cells.Item(457, 5).Value = "FILIPSTAD"
$ws.Cells.Item(457, 7).Value = 1.8
$ws.Cells.Item(457, 8).Value = 0
$ws.Cells.Item(457, 9).Value = 0
$ws.Cells.Item(457, 10).Value = 0
$ws.Cells.Item(457, 11).Value = 0
$ws.Cells.Item(457, 12).Value = 0
$ws.Cells.Item(457, 13).Value = 0
$ws.Cells.Item(457, 14).Value = 0
$ws.Cells.Item(457, 15).Value = 0
$ws.Cells.Item(457, 16).Value = 0
$ws.Cells.Item(457, 17).Value = 0
$ws.Cells.Item(457, 18).WrapText = $true
$ws.Rows.Item(457).RowHeight = 15

# 4. Add new row 458
$ws.Cells.Item(458, 1).Value = "A 46513-2023"
$ws.Cells.Item(458, 2).Value = 45197
$ws.Cells.Item(458, 2).NumberFormat = "YYYY-MM-DD"
$ws.Cells.Item(458, 3).Value = 45202
$ws.Cells.Item(458, 3).NumberFormat = "YYYY-MM-DD"
$ws.Cells.Item(458, 4).Value = "VÄRMLANDS LÄN"
$ws.Cells.Item(458, 5).Value = "FILIPSTAD"
$ws.Cells.Item(458, 6).Value = "Bergvik skog väst AB"
$ws.Cells.Item(458, 7).Value = 3.9
$ws.Cells.Item(458, 8).Value = 0
$ws.Cells.Item(458, 9).Value = 0
$ws.Cells.Item(458, 10).Value = 0
$ws.Cells.Item(458, 11).Value = 0
$ws.Cells.Item(458, 12).Value = 0
$ws.Cells.Item(458, 13).Value = 0
$ws.Cells.Item(458, 14).Value = 0
$ws.Cells.Item(458, 15).Value = 0
$ws.Cells.Item(458, 16).Value = 0
$ws.Cells.Item(458, 17).Value = 0
$ws.Cells.Item(458, 18).WrapText = $true
